$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 9.5
$ws.Range("F6").Value = 9.5
$ws.Range("G6").Value = 9

$ws.Range("E7").Value = "The complexity would be O(countries + borders*colors). You can also consider the number of colors constant if they never change, so the complexity would be just O(countries + borders)."
$ws.Range("F7").Value = "Very good. Please, use PDF version of the documents. It would make more sense to indicate the specific examples in the Test class, not in the implementation one."
$ws.Range("G7").Value = "Bactracking has a O(3^n) complexity. You got some weird values sometimes (backrtacking pruning (with balancing condition) cannot be better than backtracking without pruning in any case"

$ws.Range("A7").Select()
